# Update workbook / worksheet to reflect data refresh through 2021-10-03
# (commit message: "Add data for 2021-10-03"), which bumped the "as of"
# date referenced in the sheet name / label from 09-24 to 09-25 and
# updated the September and Total rows' 2016-2021 values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet (tab) and update its title text.
$ws.Name = "Through 2021-09-25"

# Update the "September (through 09-24)" label cell.
$ws.Range("A10").Value = "September (through 09-25)"

# September row (row 10): update 2016-2021 values (columns C-H).
$ws.Range("C10").Value = 36
$ws.Range("D10").Value = 64
$ws.Range("E10").Value = 49
$ws.Range("F10").Value = 61
$ws.Range("G10").Value = 98
$ws.Range("H10").Value = 152

# Total row (row 11): update 2016-2021 values (columns C-H).
$ws.Range("C11").Value = 417
$ws.Range("D11").Value = 615
$ws.Range("E11").Value = 539
$ws.Range("F11").Value = 410
$ws.Range("G11").Value = 882
$ws.Range("H11").Value = 1222
